$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from H1 (bold, bordered, centered) onto the new header cells
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))

# Set header text for new columns I (I0) and J (IF)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF), rows 2-19
$data = @(
    @(7, 9),
    @(7, 8),
    @(8, 9),
    @(8, 8),
    @(8, 8),
    @(7, 8),
    @(8, 9),
    @(8, 8),
    @(6, 7),
    @(7, 9),
    @(8, 9),
    @(4, 6),
    @(7, 7),
    @(7, 8),
    @(6, 7),
    @(8, 8),
    @(3, 4),
    @(4, 5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
